# Generate Report for Handoff
# Updates the localization-status report: the "Status" column moves from
# "Handed back: in sync with en-US" to "Ready for handoff" on every sheet
# that shows it, and the associated timestamp cells are refreshed to
# reflect the new handoff generation time.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status column(s): "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handoff Datetime timestamps
$overview.Range("G2").Value = "2016-08-23 06:57:03"
$dede.Range("H2").Value = "2016-08-23 06:57:03"
$zhcn.Range("H2").Value = "2016-08-23 06:56:56"

# The Status column now holds the shorter "Ready for handoff" text, so the
# sheet/table re-fits those columns to the new (narrower) content width.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
